$wb = $excel.ActiveWorkbook

# Update the IT sheet's Initial Time year value from 2020 to 2021
$itSheet = $wb.Worksheets.Item("IT")
$itSheet.Range("B2").Value = 2021

# About sheet: clear the redundant bold style duplicate from cell A7
# (removes the now-unused cellXfs entry, reverting A7 to the default style)
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("A7").Style = "Normal"

# Make the "About" sheet the active/selected tab (tabSelected moves from IT to About)
$aboutSheet.Activate()
